# ITSupplyCustomers.xlsx - "Add files via upload" edit
#
# Summary of changes applied to the "contacts" sheet (3rd sheet):
#  - Re-numbered the contactID column (A) from 1..92 to 1001..1092
#  - Fixed a zip-code typo in row 3 (H3: 924203 -> 94203)
#  - Split the malformed "DBA (larry)" contact name in row 34 into
#    separate first/last name cells (B34 = "Larry", C34 = "DBA"),
#    highlighted in yellow to flag the cleanup like other fixed rows
#  - Gave two wrapped-text rows explicit (custom) row heights so they
#    render the same regardless of default row sizing
#  - Widened columns E and I so the longer address/phone text fits
#  - Made "contacts" the active sheet/tab, with a new cell selection
#
# Sheet1 and customers sheets keep their data; only their window/view
# state (scroll position, selected tab) changes as a side effect of
# contacts becoming the active sheet.

$wb = $excel.ActiveWorkbook

$wsSheet1    = $wb.Worksheets.Item("Sheet1")
$wsCustomers = $wb.Worksheets.Item("customers")
$wsContacts  = $wb.Worksheets.Item("contacts")

# --- contacts sheet: data cleanup -----------------------------------

# Re-number contactID values in column A (rows 3-94) from 1..92 to 1001..1092
for ($r = 3; $r -le 94; $r++) {
    $oldId = $wsContacts.Cells.Item($r, 1).Value2
    $wsContacts.Cells.Item($r, 1).Value = $oldId + 1000
}

# Fix zip code typo for the first contact row
$wsContacts.Range("H3").Value = 94203

# Split "DBA (larry)" into separate last name / first name cells and
# highlight them yellow (same treatment used for other corrected rows)
$wsContacts.Range("C34").Value = "DBA"
$wsContacts.Range("B34").Value = "Larry"
$wsContacts.Range("B34:C34").Interior.Color = 65535

# Give the two wrapped multi-line rows an explicit custom row height
$wsContacts.Rows.Item(42).RowHeight = 32
$wsContacts.Rows.Item(43).RowHeight = 32
$wsContacts.Rows.Item(90).RowHeight = 32
$wsContacts.Rows.Item(91).RowHeight = 32
$wsContacts.Rows.Item(92).RowHeight = 32
$wsContacts.Rows.Item(93).RowHeight = 48
$wsContacts.Rows.Item(94).RowHeight = 48

# Widen column E (address) and column I (phone number) so values fit
$wsContacts.Columns.Item(5).ColumnWidth = 21.7
$wsContacts.Columns.Item(9).ColumnWidth = 12.8

# --- window/view state ------------------------------------------------

# Make "contacts" the active sheet/tab and move the selection
$wsContacts.Activate()
$wsContacts.Range("N7").Select()
